# Apply the "0. ตรวจ Report" refresh: new equipment rows, extended print
# area, duplicate-value highlighting on the Certificate No. column, and a
# moved selection — mirrors the authored workbook's xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ใบขอรับบริการสอบเทียบ")

# ---------------------------------------------------------------------
# 1) Replace the data rows (2-16) with the refreshed equipment list.
#    Rows 15-16 were blank before; copy row 14's cell formatting down
#    into them first so the new values land with matching borders /
#    alignment (style index 10 on column H, 1/2 on the rest).
# ---------------------------------------------------------------------
$ws.Range("A14:H14").Copy()
$ws.Range("A15:H16").PasteSpecial(-4122)  # xlPasteFormats

$rows = @(
    @{ Row=2; "A"="HEMATOCRIT CENTRIFUGE"; "B"="HETTICH"; "C"="HAEMATOKRIT 200"; "D"="0004194-04"; "E"="6515-006-0002/6/63"; "F"=45873; "G"="CF25084695"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=3; "A"="CENTRIFUGE"; "B"="PREMIERE"; "C"="XC-2008"; "D"="YDX296-3"; "E"="6515-006-0001/6/67"; "F"=45873; "G"="CF25084696"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=4; "A"="CENTRIFUGE"; "B"="JOANLAB"; "C"="LC5000H"; "D"="202401246401"; "E"="6515-006-0002/8/67"; "F"=45873; "G"="CF25084697"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=5; "A"="CENTRIFUGE"; "B"="NUVE"; "C"="NF 200"; "D"="02-8191"; "E"="6515-006-0001/2/57"; "F"=45873; "G"="CF25084698"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=6; "A"="CENTRIFUGE"; "B"="SHUKE"; "C"="TD-5M"; "D"=20040321; "E"="6515-006-0001/5/64"; "F"=45873; "G"="CF25084699"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=7; "A"="CENTRIFUGE"; "B"="CENTURION SCIENTIFIC"; "C"="C2004"; "D"="214262-11"; "E"="6515-006-0001/1/57"; "F"=45873; "G"="CF25084700"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=8; "A"="CENTRIFUGE"; "B"="BIO RAD"; "C"="DIACENT-12"; "D"=2002531; "E"="6515-006-0001/4/63"; "F"=45873; "G"="CF25084701"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=9; "A"="REFRIGERATOR"; "B"="HAIER"; "C"="HYC-360"; "D"="BE032KE1T00QEFBR0003"; "E"="6515-038-5300/1/59"; "F"=45873; "G"="CH25084702"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=10; "A"="REFRIGERATOR"; "B"="SIAMATIC"; "C"="HURR PLUS 3 DS"; "D"="HUP-03-L0265-0665-019"; "E"="6515-038-5300/2/66"; "F"=45873; "G"="CH25084703"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=11; "A"="REFRIGERATOR"; "B"="LICC"; "C"="BXC-160"; "D"="1701001909BXC-160"; "E"="6515-038-2501/2/63"; "F"=45873; "G"="CH25084704"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=12; "A"="FREEZER"; "B"="VESTFROST"; "C"="VT306"; "D"=20183227545; "E"="6515-038-3091/2/63"; "F"=45873; "G"="CH25084705"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=13; "A"="FREEZER"; "B"="HAIER"; "C"="DW-40L92"; "D"="BE02M8E0N00QEG1T0003"; "E"="6515-038-3091/1/59"; "F"=45873; "G"="CH25084706"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=14; "A"="REFRIGERATOR"; "B"="HELMER"; "C"="IBR120-GX"; "D"=2120096; "E"="6515-038-2501/3/65"; "F"=45873; "G"="CH25084707"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=15; "A"="ROTATOR"; "B"="GEMMY"; "C"="VRN-360"; "D"=1303370; "E"="6515-024-0004/1/57"; "F"=45873; "G"="CF25084708"; "H"="KOHCHAN HOSPITAL" },
    @{ Row=16; "A"="WATER BATH"; "B"="DAIHAN SCIENTIFIC"; "C"="WB-11"; "D"="0400701158C011"; "E"="6515-038-5003/2/59"; "F"=45873; "G"="WB25084709"; "H"="KOHCHAN HOSPITAL" }
)

foreach ($r in $rows) {
    $rowNum = $r["Row"]
    foreach ($col in @("A","B","C","D","E","F","G","H")) {
        $ws.Range("$col$rowNum").Value = $r[$col]
    }
}

# ---------------------------------------------------------------------
# 2) Extend the print area from A1:H14 to A1:H16.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name() -like "*Print_Area*") {
        $n.RefersTo = "=ใบขอรับบริการสอบเทียบ!`$A`$1:`$H`$16"
    }
}

# ---------------------------------------------------------------------
# 3) Add duplicate-value highlighting on column G (Certificate No.),
#    matching the existing rule already present on column E (ID No.).
# ---------------------------------------------------------------------
$rngG = $ws.Range("G1:G1048576")
$fcG = $rngG.FormatConditions.AddUniqueValues()
$fcG.DupeUnique = 1
$fcG.Font.Color = 393372
$fcG.Interior.Color = 13551615

# ---------------------------------------------------------------------
# 4) Move the active selection to M9.
# ---------------------------------------------------------------------
$ws.Range("M9").Select()
